$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the URL hyperlink (B1): new target, old display text retained
$urlCell = $ws.Range("B1")
$urlCell.Hyperlinks.Delete()
$ws.Hyperlinks.Add($urlCell, "https://www.expedia.com/", [Type]::Missing, [Type]::Missing, "https://google.com.vn")
$urlCell.Value = "https://www.expedia.com/"
$urlCell.Style = "Hyperlink"

# Row 2 - Element Name
$ws.Range("B2").Value = "tab-flight-btn-id"
$ws.Range("C2").Value = "roundtrip-btn-id"
$ws.Range("D2").Value = "flignt-origin-txt-id"
$ws.Range("E2").Value = "flight-destination-txt-id"
$ws.Range("F2").Value = "flight-add-hotel-ckb-id"
$ws.Range("G2").Value = "flight-departing-txt-id"
$ws.Range("H2").Value = "search-btn-class"

# Row 3 - Element Value
$ws.Range("B3").Value = "tab-flight-tab-hp"
$ws.Range("C3").Value = "flight-type-roundtrip-label-hp-flight"
$ws.Range("D3").Value = "flight-origin-hp-flight"
$ws.Range("E3").Value = "flight-destination-hp-flight"
$ws.Range("F3").Value = "flight-add-hotel-checkbox-hp-flight"
$ws.Range("G3").Value = "flight-departing-hp-flight"
$ws.Range("H3").Value = "gcw-submit"

$ws.Range("C1").Style = "Hyperlink"

Write-Host "done"
